$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A1").Value = "Derivados del petróleo"
$ws.Range("A2").Value = "Madera"
$ws.Range("A3").Value = " No aplicable"
$ws.Range("A4").Value = "Otros"
$ws.Range("A5").Value = "Gas"
$ws.Range("A6").Value = "Electricidad"
$ws.Range("A7").Value = "Derivados del carbón"
